# cryptos.xlsx refresh - "Updated cryptos list ... with GitHub Actions"
# Applies the scraped price / Volume(1h) updates, including the two pairs of
# coin rows that swapped rank position (rows 26/27: Cosmos<->Toncoin, rows
# 49/50: Celestia<->NEARProtocol - name/link/price/volume all move together).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) and Volume(1h) (column E) are plain text in this sheet.
# Assigning a numeric-looking string straight to .Value makes Excel silently
# reinterpret it as a real number (dropping significant trailing/leading
# zeros such as "230.50" or "0.0000100"), so those writes go through
# Set-TextValue, which quote-prefixes the string first - the same thing a
# user typing  '230.51  into a cell does - to force it to stay text.
function Set-TextValue($addr, $value) {
    $ws.Range($addr).Value = "'" + $value
}
# Row 2
$ws.Range('D2').Value = '43.995.25'
$ws.Range('E2').Value = '  +0.42%  '

# Row 3
$ws.Range('D3').Value = '2.259.74'
$ws.Range('E3').Value = '  -0.28%  '

# Row 4
$ws.Range('E4').Value = '  +0.28%  '

# Row 5
Set-TextValue 'D5' '230.51'
$ws.Range('E5').Value = '  -0.10%  '

# Row 6
Set-TextValue 'D6' '0.639'
$ws.Range('E6').Value = '  +1.95%  '

# Row 7
Set-TextValue 'D7' '64.22'
$ws.Range('E7').Value = '  +5.05%  '

# Row 8
$ws.Range('E8').Value = '  +0.10%  '

# Row 9
Set-TextValue 'D9' '0.450'
$ws.Range('E9').Value = '  +7.06%  '

# Row 10
Set-TextValue 'D10' '0.0992'
$ws.Range('E10').Value = '  +6.27%  '

# Row 11
Set-TextValue 'D11' '57.10'
$ws.Range('E11').Value = '  -1.61%  '

# Row 12
$ws.Range('E12').Value = '  +14.79%  '

# Row 13
$ws.Range('E13').Value = '  +1.91%  '

# Row 14
$ws.Range('D14').Value = '2.596.57'
$ws.Range('E14').Value = '  -0.31%  '

# Row 15
Set-TextValue 'D15' '15.69'
$ws.Range('E15').Value = '  +0.69%  '

# Row 16
$ws.Range('E16').Value = '  +4.80%  '

# Row 17
Set-TextValue 'D17' '0.833'
$ws.Range('E17').Value = '  +2.92%  '

# Row 18
$ws.Range('D18').Value = '2.261.26'
$ws.Range('E18').Value = '  -0.22%  '

# Row 19
$ws.Range('D19').Value = '43.890.67'
$ws.Range('E19').Value = '  +0.74%  '

# Row 20
Set-TextValue 'D20' '0.0000100'
$ws.Range('E20').Value = '  +7.17%  '

# Row 21
Set-TextValue 'D21' '73.38'
$ws.Range('E21').Value = '  +0.70%  '

# Row 22
Set-TextValue 'D22' '6.06'
$ws.Range('E22').Value = '  -2.41%  '

# Row 23
Set-TextValue 'D23' '250.81'
$ws.Range('E23').Value = '  -1.13%  '

# Row 24
$ws.Range('E24').Value = '  -0.02%  '

# Row 25
$ws.Range('E25').Value = '  -4.16%  '

# Row 26
$ws.Range('B26').Value = 'Toncoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 'D26' '2.30'
$ws.Range('E26').Value = '  +0.51%  '

# Row 27
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue 'D27' '10.08'
$ws.Range('E27').Value = '  +2.40%  '

# Row 28
Set-TextValue 'D28' '3.26'
$ws.Range('E28').Value = '  +22.59%  '

# Row 29
Set-TextValue 'D29' '170.87'
$ws.Range('E29').Value = '  +0.22%  '

# Row 30
Set-TextValue 'D30' '0.139'
$ws.Range('E30').Value = '  -0.37%  '

# Row 31
Set-TextValue 'D31' '20.90'
$ws.Range('E31').Value = '  +2.04%  '

# Row 32
$ws.Range('E32').Value = '  -3.58%  '

# Row 33
Set-TextValue 'D33' '0.125'
$ws.Range('E33').Value = '  +2.57%  '

# Row 34
Set-TextValue 'D34' '0.0703'
$ws.Range('E34').Value = '  +6.64%  '

# Row 35
Set-TextValue 'D35' '4.77'
$ws.Range('E35').Value = '  -0.24%  '

# Row 36
Set-TextValue 'D36' '4.89'
$ws.Range('E36').Value = '  -3.14%  '

# Row 37
$ws.Range('E37').Value = '  +5.09%  '

# Row 38
Set-TextValue 'D38' '6.47'
$ws.Range('E38').Value = '  +0.10%  '

# Row 39
$ws.Range('E39').Value = '  -4.09%  '

# Row 40
$ws.Range('E40').Value = '  +3.53%  '

# Row 41
$ws.Range('E41').Value = '  -0.02%  '

# Row 42
Set-TextValue 'D42' '0.000224'
$ws.Range('E42').Value = '  -1.48%  '

# Row 43
$ws.Range('E43').Value = '  -1.36%  '

# Row 44
$ws.Range('E44').Value = '  +4.16%  '

# Row 45
Set-TextValue 'D45' '8.21'
$ws.Range('E45').Value = '  -5.82%  '

# Row 46
Set-TextValue 'D46' '97.88'
$ws.Range('E46').Value = '  -0.41%  '

# Row 47
$ws.Range('E47').Value = '  -0.71%  '

# Row 48
Set-TextValue 'D48' '4.40'
$ws.Range('E48').Value = '  -2.72%  '

# Row 49
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 'D49' '2.37'
$ws.Range('E49').Value = '  +5.13%  '

# Row 50
$ws.Range('B50').Value = 'Celestia'
$ws.Range('C50').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextValue 'D50' '10.20'
$ws.Range('E50').Value = '  +6.81%  '

# Row 51
$ws.Range('D51').Value = '1.437.96'
$ws.Range('E51').Value = '  -2.30%  '
